$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    12  = 4.951000000000001
    27  = 5.516
    32  = 6.77
    36  = 8.599
    38  = 6.027
    46  = 6.311
    54  = 5.176
    55  = 4.763
    56  = 4.654000000000001
    67  = 5.516
    69  = 5.367
    72  = 5.7
    83  = 5.883
    86  = 5.02
    91  = 5.923
    93  = 5.424000000000001
    99  = 5.699
    104 = 7.950999999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}

$wb.Save()
